$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number and date range) ---
$ws.Range("A8").Value = "Volume 30   Number  33"
$ws.Range("C9").Value = "Report Covering the Week  8/14/2023  Through  8/20/2023"

# --- Plain numeric value updates (no type/style change) ---
$ws.Range("L15").Value = -22.222222222222
$ws.Range("M15").Value = -26.315789473684
$ws.Range("N15").Value = -26.315789473684
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 32
$ws.Range("H16").Value = 18.518518518518
$ws.Range("I16").Value = 203
$ws.Range("J16").Value = 176
$ws.Range("K16").Value = 15.340909090909
$ws.Range("L16").Value = 34.437086092715
$ws.Range("M16").Value = -5.140186915887
$ws.Range("N16").Value = -75.804529201430
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 80
$ws.Range("F17").Value = 41
$ws.Range("G17").Value = 24
$ws.Range("H17").Value = 70.833333333333
$ws.Range("I17").Value = 300
$ws.Range("J17").Value = 277
$ws.Range("K17").Value = 8.303249097472
$ws.Range("L17").Value = 25
$ws.Range("M17").Value = 44.230769230769
$ws.Range("N17").Value = 11.111111111111
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -66.666666666666
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 102
$ws.Range("J18").Value = 96
$ws.Range("K18").Value = 6.25
$ws.Range("L18").Value = 25.925925925925
$ws.Range("M18").Value = -46.031746031746
$ws.Range("N18").Value = -92.296072507552
$ws.Range("C19").Value = 17
$ws.Range("D19").Value = 17
$ws.Range("E19").Value = 0
$ws.Range("G19").Value = 90
$ws.Range("H19").Value = -31.111111111111
$ws.Range("I19").Value = 524
$ws.Range("J19").Value = 623
$ws.Range("K19").Value = -15.890850722311
$ws.Range("L19").Value = 59.756097560975
$ws.Range("M19").Value = 73.509933774834
$ws.Range("N19").Value = -45.015739769150
$ws.Range("C20").Value = 8
$ws.Range("E20").Value = -11.111111111111
$ws.Range("G20").Value = 25
$ws.Range("H20").Value = 4
$ws.Range("I20").Value = 199
$ws.Range("J20").Value = 191
$ws.Range("K20").Value = 4.188481675392
$ws.Range("L20").Value = 54.263565891472
$ws.Range("M20").Value = 37.241379310344
$ws.Range("N20").Value = -85.966149506347
$ws.Range("C21").Value = 42
$ws.Range("D21").Value = 41
$ws.Range("E21").Value = 2.439024390243
$ws.Range("F21").Value = 175
$ws.Range("G21").Value = 182
$ws.Range("H21").Value = -3.846153846153
$ws.Range("I21").Value = 1344
$ws.Range("J21").Value = 1389
$ws.Range("K21").Value = -3.239740820734
$ws.Range("L21").Value = 41.324921135646
$ws.Range("M21").Value = 24.675324675324
$ws.Range("N21").Value = -72.208436724565
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = -16.666666666666
$ws.Range("I22").Value = 61
$ws.Range("K22").Value = 74.285714285714
$ws.Range("L22").Value = 281.25
$ws.Range("M22").Value = 165.217391304348
$ws.Range("C24").Value = 33
$ws.Range("D24").Value = 36
$ws.Range("E24").Value = -8.333333333333
$ws.Range("F24").Value = 139
$ws.Range("G24").Value = 198
$ws.Range("H24").Value = -29.797979797979
$ws.Range("I24").Value = 1274
$ws.Range("J24").Value = 1188
$ws.Range("K24").Value = 7.239057239057
$ws.Range("L24").Value = 49.180327868852
$ws.Range("M24").Value = 82
$ws.Range("D25").Value = 16
$ws.Range("E25").Value = 25
$ws.Range("F25").Value = 77
$ws.Range("G25").Value = 71
$ws.Range("H25").Value = 8.450704225352
$ws.Range("I25").Value = 590
$ws.Range("J25").Value = 576
$ws.Range("K25").Value = 2.430555555555
$ws.Range("L25").Value = 11.742424242424
$ws.Range("M25").Value = -1.830282861896
$ws.Range("E26").Value = -100
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = -40
$ws.Range("J26").Value = 32
$ws.Range("K26").Value = -6.25
$ws.Range("L26").Value = -9.090909090909
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = 87.5
$ws.Range("I27").Value = 102
$ws.Range("J27").Value = 65
$ws.Range("K27").Value = 56.923076923076
$ws.Range("L27").Value = 43.661971830985

# --- Cells changing from numeric to text placeholder ("0" or "***.*") ---
# Row 14 is never touched by this edit, so C14 (style14,"0") and E14 (style14,"***.*")
# are used as stable format-source cells.
$ws.Range("D15").Value = "'0"
$ws.Range("E15").Value = "'***.*"
$ws.Range("D22").Value = "'0"
$ws.Range("E22").Value = "'***.*"
$ws.Range("C26").Value = "'0"
$ws.Range("F28").Value = "'0"
$ws.Range("F29").Value = "'0"
$ws.Range("D30").Value = "'0"
$ws.Range("E30").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("F28").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E30").PasteSpecial(-4122)

# --- Cells changing from text placeholder to numeric values ---
# Row 14 I14 (style15, integer columns) / K14 (style16, %-chg columns) are stable format sources.
# The correct source is determined by which column the cell is in (count columns use style 15,
# %-change columns use style 16), NOT by whether the target value happens to contain a decimal point.
$ws.Range("C22").Value = 3
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 200
$ws.Range("I14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("I14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("K14").Copy()
$ws.Range("E27").PasteSpecial(-4122)

$excel.CutCopyMode = $false